$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Add the ORCID userId for Bethany Fowler (row 7, column F = userId)
$ws.Range("F7").Value = "0000-0001-8655-7253"

# Update the active selection to match the edited cell
$ws.Activate()
$ws.Range("F7").Select()
